# Update the "K" (strikeouts) column (G) values on the active sheet.
# The raw Strike# values previously stored in column G are being
# regenerated/replaced with the recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 9
    3  = 2
    4  = 3
    5  = 3
    6  = 6
    7  = 5
    8  = 3
    9  = 9
    10 = 7
    11 = 11
    12 = 6
    13 = 6
    14 = 7
    15 = 6
    16 = 6
    17 = 1
    18 = 4
    19 = 3
    20 = 4
    21 = 2
    22 = 2
    23 = 4
    24 = 4
    25 = 5
    26 = 3
    27 = 6
    28 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
